# Weekly update: insert a new "Mango" price record for Vega Monumental
# Concepción at row 123, shifting the existing historical rows (123-183)
# down by one row (124-184). The new row carries this week's figures;
# every other row keeps the same category metadata and simply slides
# down to make room, exactly as Excel does on Rows.Insert().

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 123:183 down to 124:184, leaving a fresh blank row 123.
$ws.Rows.Item(123).Insert()

# Populate the new row 123 with this week's observation.
$ws.Range("A123").Value = 11
$ws.Range("B123").Value = "Vega Monumental Concepción"
$ws.Range("C123").Value = "Bíobío"
$ws.Range("D123").Value = 45134
$ws.Range("E123").Value = 8
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100108
$ws.Range("H123").Value = "Tropicales y subtropicales"
$ws.Range("I123").Value = 100108002
$ws.Range("J123").Value = "Mango"
$ws.Range("K123").Value = "Sin especificar"
$ws.Range("L123").Value = "Primera"
$ws.Range("M123").Value = 150
$ws.Range("N123").Value = 8000
$ws.Range("O123").Value = 8000
$ws.Range("P123").Value = 8000
$ws.Range("Q123").Value = "`$/bandeja 4 kilos"
$ws.Range("R123").Value = "Perú"
$ws.Range("S123").Value = 2000
$ws.Range("T123").Value = 4
